$wb = $excel.ActiveWorkbook

# 1) Rename header cells on existing sheets (Requested quantity -> *_PO_Qty)
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2) Add the new "PO Forecast" sheet after the last existing sheet
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# 3) Header row values
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy the bold/centered/bordered header formatting used on the other sheets
$srcHeader = $wsWeekly.Range("A1:B1")
$srcHeader.Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$wsForecast.Range("A2").Value = 45389.99999999999
$wsForecast.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B2").Value = 59
$wsForecast.Range("C2").Value = -44.44967783979528
$wsForecast.Range("D2").Value = 167.9470193144836

$wsForecast.Range("A3").Value = 45396.99999999999
$wsForecast.Range("A3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B3").Value = 60
$wsForecast.Range("C3").Value = -51.08012032099543
$wsForecast.Range("D3").Value = 172.5043106071817

$wsForecast.Range("A4").Value = 45403.99999999999
$wsForecast.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B4").Value = 60
$wsForecast.Range("C4").Value = -59.71832749387451
$wsForecast.Range("D4").Value = 173.9362087666918

$wsForecast.Range("A5").Value = 45410.99999999999
$wsForecast.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B5").Value = 60
$wsForecast.Range("C5").Value = -55.87343863609873
$wsForecast.Range("D5").Value = 173.4923630946206

$wsForecast.Range("A6").Value = 45417.99999999999
$wsForecast.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B6").Value = 60
$wsForecast.Range("C6").Value = -55.41673956397967
$wsForecast.Range("D6").Value = 173.7389176877349

$wsForecast.Range("A7").Value = 45424.99999999999
$wsForecast.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B7").Value = 61
$wsForecast.Range("C7").Value = -43.6413660606327
$wsForecast.Range("D7").Value = 177.4904579356019

$wsForecast.Range("A8").Value = 45431.99999999999
$wsForecast.Range("A8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B8").Value = 61
$wsForecast.Range("C8").Value = -42.65928442242895
$wsForecast.Range("D8").Value = 167.2477957495742

$wsForecast.Range("A9").Value = 45438.99999999999
$wsForecast.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B9").Value = 61
$wsForecast.Range("C9").Value = -54.9895491999991
$wsForecast.Range("D9").Value = 170.6493389902022

$wsForecast.Range("A10").Value = 45445.99999999999
$wsForecast.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B10").Value = 61
$wsForecast.Range("C10").Value = -52.59389791907213
$wsForecast.Range("D10").Value = 166.7940210216497

$wsForecast.Range("A11").Value = 45452.99999999999
$wsForecast.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B11").Value = 61
$wsForecast.Range("C11").Value = -46.49175250478561
$wsForecast.Range("D11").Value = 163.8513673536992

$wsForecast.Range("A12").Value = 45459.99999999999
$wsForecast.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B12").Value = 62
$wsForecast.Range("C12").Value = -54.5073009215438
$wsForecast.Range("D12").Value = 179.5223813679918

$wsForecast.Range("A13").Value = 45466.99999999999
$wsForecast.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B13").Value = 62
$wsForecast.Range("C13").Value = -51.89812795714681
$wsForecast.Range("D13").Value = 167.4499926266766

$wsForecast.Range("A14").Value = 45473.99999999999
$wsForecast.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B14").Value = 62
$wsForecast.Range("C14").Value = -49.98822424468141
$wsForecast.Range("D14").Value = 174.5396328595442

$wsForecast.Range("A15").Value = 45480.99999999999
$wsForecast.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B15").Value = 62
$wsForecast.Range("C15").Value = -45.20979875988159
$wsForecast.Range("D15").Value = 178.2225860524472

$wsForecast.Range("A16").Value = 45487.99999999999
$wsForecast.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B16").Value = 63
$wsForecast.Range("C16").Value = -53.54288427415901
$wsForecast.Range("D16").Value = 175.3035008310999

$wsForecast.Range("A17").Value = 45494.99999999999
$wsForecast.Range("A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B17").Value = 63
$wsForecast.Range("C17").Value = -44.1289494554
$wsForecast.Range("D17").Value = 185.2617533061213

$wsForecast.Range("A18").Value = 45501.99999999999
$wsForecast.Range("A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B18").Value = 63
$wsForecast.Range("C18").Value = -47.02459564306704
$wsForecast.Range("D18").Value = 181.3452051746392

$wsForecast.Range("A19").Value = 45508.99999999999
$wsForecast.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B19").Value = 63
$wsForecast.Range("C19").Value = -48.48097108541676
$wsForecast.Range("D19").Value = 178.9000421259577

$wsForecast.Range("A20").Value = 45515.99999999999
$wsForecast.Range("A20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B20").Value = 63
$wsForecast.Range("C20").Value = -50.79077697361235
$wsForecast.Range("D20").Value = 170.1844913544595

$wsForecast.Range("A21").Value = 45522.99999999999
$wsForecast.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B21").Value = 64
$wsForecast.Range("C21").Value = -50.57811785219108
$wsForecast.Range("D21").Value = 183.3697555759855

$wsForecast.Range("A22").Value = 45529.99999999999
$wsForecast.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B22").Value = 64
$wsForecast.Range("C22").Value = -53.60452573837925
$wsForecast.Range("D22").Value = 179.1574771931968

$wsForecast.Range("A23").Value = 45536.99999999999
$wsForecast.Range("A23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B23").Value = 64
$wsForecast.Range("C23").Value = -48.2677292921198
$wsForecast.Range("D23").Value = 175.5048637322905

$wsForecast.Range("A24").Value = 45571.99999999999
$wsForecast.Range("A24").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B24").Value = 65
$wsForecast.Range("C24").Value = -46.08669165559618
$wsForecast.Range("D24").Value = 170.9225537191958

$wsForecast.Range("A25").Value = 45578.99999999999
$wsForecast.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B25").Value = 66
$wsForecast.Range("C25").Value = -46.36547736957442
$wsForecast.Range("D25").Value = 176.7035621448274

$wsForecast.Range("A26").Value = 45585.99999999999
$wsForecast.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B26").Value = 66
$wsForecast.Range("C26").Value = -40.73406853691953
$wsForecast.Range("D26").Value = 172.1690761991554

$wsForecast.Range("A27").Value = 45606.99999999999
$wsForecast.Range("A27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B27").Value = 66
$wsForecast.Range("C27").Value = -48.32186878909408
$wsForecast.Range("D27").Value = 173.1479675667509

$wsForecast.Range("A28").Value = 45613.99999999999
$wsForecast.Range("A28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B28").Value = 67
$wsForecast.Range("C28").Value = -47.22924310069542
$wsForecast.Range("D28").Value = 180.267166366457

$wsForecast.Range("A29").Value = 45620.99999999999
$wsForecast.Range("A29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B29").Value = 67
$wsForecast.Range("C29").Value = -47.19759372068311
$wsForecast.Range("D29").Value = 184.9119848835122

$wsForecast.Range("A30").Value = 45627.99999999999
$wsForecast.Range("A30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B30").Value = 67
$wsForecast.Range("C30").Value = -47.10551211782278
$wsForecast.Range("D30").Value = 172.3438260318052

$wsForecast.Range("A31").Value = 45634.99999999999
$wsForecast.Range("A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B31").Value = 67
$wsForecast.Range("C31").Value = -49.32748982426994
$wsForecast.Range("D31").Value = 181.1089250232354

$wsForecast.Range("A32").Value = 45641.99999999999
$wsForecast.Range("A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B32").Value = 68
$wsForecast.Range("C32").Value = -49.31567159947596
$wsForecast.Range("D32").Value = 182.529928116671

$wsForecast.Range("A33").Value = 45648.99999999999
$wsForecast.Range("A33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B33").Value = 68
$wsForecast.Range("C33").Value = -46.73163403028917
$wsForecast.Range("D33").Value = 175.7505482252068

$wsForecast.Range("A34").Value = 45655.99999999999
$wsForecast.Range("A34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B34").Value = 68
$wsForecast.Range("C34").Value = -41.72391015590549
$wsForecast.Range("D34").Value = 183.8593638430806

$wsForecast.Range("A35").Value = 45662.99999999999
$wsForecast.Range("A35").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B35").Value = 68
$wsForecast.Range("C35").Value = -37.35975068610952
$wsForecast.Range("D35").Value = 181.8294974001026

$wsForecast.Range("A36").Value = 45669.99999999999
$wsForecast.Range("A36").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B36").Value = 68
$wsForecast.Range("C36").Value = -37.84546035434995
$wsForecast.Range("D36").Value = 173.8157977740092

$wsForecast.Range("A37").Value = 45676.99999999999
$wsForecast.Range("A37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B37").Value = 69
$wsForecast.Range("C37").Value = -36.93226911486622
$wsForecast.Range("D37").Value = 179.8096396754676

$wsForecast.Range("A38").Value = 45683.99999999999
$wsForecast.Range("A38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B38").Value = 69
$wsForecast.Range("C38").Value = -44.58494590714573
$wsForecast.Range("D38").Value = 182.3095754989001

$wsForecast.Range("A39").Value = 45690.99999999999
$wsForecast.Range("A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B39").Value = 69
$wsForecast.Range("C39").Value = -33.91870549064905
$wsForecast.Range("D39").Value = 190.1240938980883

$wsForecast.Range("A40").Value = 45697.99999999999
$wsForecast.Range("A40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B40").Value = 69
$wsForecast.Range("C40").Value = -41.08600714248907
$wsForecast.Range("D40").Value = 182.1676380252661

